# Fix departure/arrival times on several rows of the "kniha jazd" (log book)
# sheet. Times are stored as plain text (inlineStr), so we explicitly set
# NumberFormat to "@" (text) before assigning the values to avoid Excel
# auto-converting "14:00" style strings into time serial values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Odchod Čas 15:00 -> 14:00, Príchod Čas 16:00 -> 15:00
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "14:00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "15:00"

# Row 5: Odchod Čas 16:00 -> 15:45, Príchod Čas 17:00 -> 16:45
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "15:45"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "16:45"

# Row 7: Odchod Čas 15:30 -> 15:15, Príchod Čas 16:30 -> 16:15
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "15:15"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "16:15"

# Row 15: Odchod Čas 16:30 -> 16:00, Príchod Čas 17:30 -> 17:00
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "16:00"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "17:00"
